$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; existing rows 45-63 shift down to 46-64.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price record.
$ws.Cells.Item(45, 1).Value  = 4
$ws.Cells.Item(45, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value  = "Los Lagos"
$ws.Cells.Item(45, 4).Value  = 44981
$ws.Cells.Item(45, 5).Value  = 10
$ws.Cells.Item(45, 6).Value  = "Fruta"
$ws.Cells.Item(45, 7).Value  = 100101
$ws.Cells.Item(45, 8).Value  = "Berries"
$ws.Cells.Item(45, 9).Value  = 100101001
$ws.Cells.Item(45, 10).Value = "Arándano (blue)"
$ws.Cells.Item(45, 11).Value = "Sin especificar"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 200
$ws.Cells.Item(45, 14).Value = 2000
$ws.Cells.Item(45, 15).Value = 2200
$ws.Cells.Item(45, 16).Value = 2100
$ws.Cells.Item(45, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(45, 19).Value = 1050
$ws.Cells.Item(45, 20).Value = 2
